$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2:D51").NumberFormat = "@"

$ws.Cells.Item(2, 4).Value = '47.110.17'
$ws.Cells.Item(2, 5).Value = '  +1.45%  '
$ws.Cells.Item(3, 4).Value = '2.487.34'
$ws.Cells.Item(3, 5).Value = '  +0.62%  '
$ws.Cells.Item(4, 5).Value = '  +0.15%  '
$ws.Cells.Item(5, 4).Value = '322.23'
$ws.Cells.Item(5, 5).Value = '  -0.29%  '
$ws.Cells.Item(6, 4).Value = '108.36'
$ws.Cells.Item(6, 5).Value = '  +2.34%  '
$ws.Cells.Item(7, 4).Value = '0.524'
$ws.Cells.Item(7, 5).Value = '  +0.89%  '
$ws.Cells.Item(8, 5).Value = '  -0.04%  '
$ws.Cells.Item(9, 4).Value = '0.533'
$ws.Cells.Item(9, 5).Value = '  -1.32%  '
$ws.Cells.Item(10, 4).Value = '38.81'
$ws.Cells.Item(10, 5).Value = '  +7.16%  '
$ws.Cells.Item(11, 4).Value = '0.0811'
$ws.Cells.Item(11, 5).Value = '  -0.71%  '
$ws.Cells.Item(12, 5).Value = '  +0.42%  '
$ws.Cells.Item(13, 4).Value = '18.31'
$ws.Cells.Item(13, 5).Value = '  -0.73%  '
$ws.Cells.Item(14, 4).Value = '7.16'
$ws.Cells.Item(14, 5).Value = '  +0.83%  '
$ws.Cells.Item(15, 4).Value = '2.877.55'
$ws.Cells.Item(15, 5).Value = '  +0.46%  '
$ws.Cells.Item(16, 4).Value = '2.488.86'
$ws.Cells.Item(16, 5).Value = '  +0.79%  '
$ws.Cells.Item(17, 4).Value = '0.848'
$ws.Cells.Item(17, 5).Value = '  +0.13%  '
$ws.Cells.Item(18, 4).Value = '47.028.36'
$ws.Cells.Item(18, 5).Value = '  +1.52%  '
$ws.Cells.Item(19, 4).Value = '12.64'
$ws.Cells.Item(19, 5).Value = '  -0.44%  '
$ws.Cells.Item(20, 4).Value = '6.58'
$ws.Cells.Item(20, 5).Value = '  +1.54%  '
$ws.Cells.Item(21, 4).Value = '0.0₃0935'
$ws.Cells.Item(21, 5).Value = '  -0.30%  '
$ws.Cells.Item(22, 4).Value = '2.69'
$ws.Cells.Item(22, 5).Value = '  +12.25%  '
$ws.Cells.Item(23, 4).Value = '70.62'
$ws.Cells.Item(23, 5).Value = '  -0.01%  '
$ws.Cells.Item(24, 4).Value = '246.37'
$ws.Cells.Item(24, 5).Value = '  -1.05%  '
$ws.Cells.Item(25, 4).Value = '2.58'
$ws.Cells.Item(25, 5).Value = '  +1.62%  '
$ws.Cells.Item(26, 5).Value = '  -0.07%  '
$ws.Cells.Item(27, 4).Value = '25.78'
$ws.Cells.Item(27, 5).Value = '  -1.69%  '
$ws.Cells.Item(28, 5).Value = '  +4.13%  '
$ws.Cells.Item(29, 4).Value = '10.01'
$ws.Cells.Item(29, 5).Value = '  +2.04%  '
$ws.Cells.Item(30, 4).Value = '0.139'
$ws.Cells.Item(30, 5).Value = '  +8.17%  '
$ws.Cells.Item(31, 4).Value = '34.99'
$ws.Cells.Item(31, 5).Value = '  +0.72%  '
$ws.Cells.Item(32, 5).Value = '  +0.73%  '
$ws.Cells.Item(33, 4).Value = '19.85'
$ws.Cells.Item(33, 5).Value = '  -0.09%  '
$ws.Cells.Item(34, 4).Value = '5.38'
$ws.Cells.Item(34, 5).Value = '  +0.55%  '
$ws.Cells.Item(35, 4).Value = '0.0783'
$ws.Cells.Item(35, 5).Value = '  +1.97%  '
$ws.Cells.Item(36, 5).Value = '  +0.30%  '
$ws.Cells.Item(37, 4).Value = '1.96'
$ws.Cells.Item(37, 5).Value = '  +2.29%  '
$ws.Cells.Item(38, 4).Value = '4.66'
$ws.Cells.Item(38, 5).Value = '  +1.24%  '
$ws.Cells.Item(39, 4).Value = '2.97'
$ws.Cells.Item(39, 5).Value = '  -0.02%  '
$ws.Cells.Item(40, 5).Value = '  +0.29%  '
$ws.Cells.Item(41, 4).Value = '120.62'
$ws.Cells.Item(41, 5).Value = '  -2.59%  '
$ws.Cells.Item(42, 5).Value = '  -1.05%  '
$ws.Cells.Item(43, 4).Value = '21.12'
$ws.Cells.Item(43, 5).Value = '  +0.93%  '
$ws.Cells.Item(44, 4).Value = '0.0295'
$ws.Cells.Item(44, 5).Value = '  +0.48%  '
$ws.Cells.Item(45, 4).Value = '1.990.15'
$ws.Cells.Item(45, 5).Value = '  +0.46%  '
$ws.Cells.Item(46, 4).Value = '3.03'
$ws.Cells.Item(46, 5).Value = '  +1.31%  '
$ws.Cells.Item(47, 4).Value = '2.04'
$ws.Cells.Item(47, 5).Value = '  -3.05%  '
$ws.Cells.Item(48, 4).Value = '1.78'
$ws.Cells.Item(48, 5).Value = '  -3.55%  '
$ws.Cells.Item(49, 4).Value = '9.07'
$ws.Cells.Item(49, 5).Value = '  -0.50%  '
$ws.Cells.Item(50, 4).Value = '5.17'
$ws.Cells.Item(50, 5).Value = '  +0.66%  '
$ws.Cells.Item(51, 4).Value = '56.20'
$ws.Cells.Item(51, 5).Value = '  +2.42%  '
